$p = $ppt.ActivePresentation

# --- 1. Table style change on the table in Slide 5 ---
$s = $p.Slides.Item(5)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{C9B8D954-E415-4624-BDCD-7A41CEB33342}")
    }
}

# --- 2. Swap the theme colour scheme from "Red Violet" (Integral) to "Office" ---
$tcs = $p.SlideMaster.Theme.ThemeColorScheme
$tcs.Item(1).RGB  = 0         # dk1     000000
$tcs.Item(2).RGB  = 16777215  # lt1     FFFFFF
$tcs.Item(3).RGB  = 6968388   # dk2     44546A
$tcs.Item(4).RGB  = 15132391  # lt2     E7E6E6
$tcs.Item(5).RGB  = 13998939  # accent1 5B9BD5
$tcs.Item(6).RGB  = 3243501   # accent2 ED7D31
$tcs.Item(7).RGB  = 10855845  # accent3 A5A5A5
$tcs.Item(8).RGB  = 49407     # accent4 FFC000
$tcs.Item(9).RGB  = 12874308  # accent5 4472C4
$tcs.Item(10).RGB = 4697456   # accent6 70AD47
$tcs.Item(11).RGB = 12673797  # hlink   0563C1
$tcs.Item(12).RGB = 7491477   # folHlink 954F72
